$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.020.61"
$ws.Range("E2").Value = "  +1.40%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.978.80"
$ws.Range("E3").Value = "  +1.08%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.77"
$ws.Range("E5").Value = "  +0.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.628"
$ws.Range("E6").Value = "  +2.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "61.02"
$ws.Range("E7").Value = "  +3.92%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0797"
$ws.Range("E10").Value = "  -1.65%  "
$ws.Range("E11").Value = "  +0.09%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "14.92"
$ws.Range("E12").Value = "  +8.87%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.28"
$ws.Range("E13").Value = "  +0.64%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.840"
$ws.Range("E14").Value = "  +1.36%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.269.47"
$ws.Range("E15").Value = "  +1.11%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.45"
$ws.Range("E16").Value = "  +3.60%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.980.46"
$ws.Range("E17").Value = "  +1.33%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "36.904.60"
$ws.Range("E18").Value = "  +1.28%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "70.05"
$ws.Range("E19").Value = "  +0.55%  "
$ws.Range("E20").Value = "  +0.31%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.15"
$ws.Range("E21").Value = "  +1.95%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "230.09"
$ws.Range("E22").Value = "  +0.80%  "
$ws.Range("E23").Value = "  +0.13%  "
$ws.Range("E24").Value = "  +2.37%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.36"
$ws.Range("E25").Value = "  +0.36%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.148"
$ws.Range("E26").Value = "  +7.66%  "
$ws.Range("E27").Value = "  +0.36%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "163.14"
$ws.Range("E28").Value = "  +1.77%  "
$ws.Range("E29").Value = "  +0.88%  "
$ws.Range("E30").Value = "  +16.91%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.122"
$ws.Range("E31").Value = "  +1.67%  "
$ws.Range("E32").Value = "  +2.82%  "
$ws.Range("E33").Value = "  +0.28%  "
$ws.Range("E34").Value = "  +5.62%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.30"
$ws.Range("E35").Value = "  +2.33%  "
$ws.Range("E36").Value = "  -0.10%  "
$ws.Range("E37").Value = "  +0.79%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.34"
$ws.Range("E38").Value = "  -0.45%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.47"
$ws.Range("E39").Value = "  -4.86%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0977"
$ws.Range("E40").Value = "  -0.17%  "
$ws.Range("E41").Value = "  +1.32%  "
$ws.Range("E42").Value = "  -0.15%  "
$ws.Range("E43").Value = "  +0.63%  "
$ws.Range("E44").Value = "  +3.59%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.369.84"
$ws.Range("E45").Value = "  +0.31%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "89.87"
$ws.Range("E46").Value = "  +2.27%  "
$ws.Range("E47").Value = "  +0.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.19"
$ws.Range("E48").Value = "  +0.91%  "
$ws.Range("E49").Value = "  -0.26%  "
$ws.Range("E50").Value = "  +6.00%  "
$ws.Range("E51").Value = "  +10.22%  "
